$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = "30 TL - 30 TL"

# Row 3
$ws.Range("F3").Value = "30,46 TL - 60,94 TL - 609,43 TL"

# Row 4
$ws.Range("F4").Value = "30,46 TL - 60,94 TL - 609,43 TL"

# Row 5
$ws.Range("F5").Value = "30,46 TL - 60,94 TL - 609,43 TL"

# Row 8
$ws.Range("F8").Value = "15,23 TL - 30,47 TL - 304,71 TL"

# Row 9
$ws.Range("F9").Value = "15,23 TL - 30,47 TL - 304,71 TL"

# Row 10
$ws.Range("F10").Value = "15,23 TL - 30,47 TL - 304,71 TL"

# Row 13
$ws.Range("C13").Value = "Hesaba: Asgari 0 TL | Azami 0,94 TL"
$ws.Range("F13").Value = "Hesaba: Asgari 300 TL | Azami 3.080 TL"

# Row 14
$ws.Range("F14").Value = "1.952,38 TL - 9.523,81 TL"
